# Apply scheduled market-data refresh to leve profit calculations across all sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 77.25
$ws.Range("I8").Value = 77.25
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 231.75
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -92.75
$ws.Range("N8").ClearContents()
$ws.Range("H19").Value = 1218.1154
$ws.Range("I19").Value = 1179.7142
$ws.Range("J19").Value = 1262.9166
$ws.Range("K19").Value = 1179.7142
$ws.Range("L19").Value = 1262.9166
$ws.Range("M19").Value = -1004.7142
$ws.Range("N19").Value = -1612.9166
$ws.Range("H41").Value = 6400
$ws.Range("I41").Value = 4333.3335
$ws.Range("J41").Value = 9500
$ws.Range("K41").Value = 4333.3335
$ws.Range("L41").Value = 9500
$ws.Range("M41").Value = -3893.3335
$ws.Range("N41").Value = -10380
$ws.Range("H62").Value = 4349
$ws.Range("I62").Value = 1499.8
$ws.Range("K62").Value = 1499.8
$ws.Range("M62").Value = -875.8
$ws.Range("H65").Value = 4349
$ws.Range("I65").Value = 1499.8
$ws.Range("K65").Value = 7499
$ws.Range("M65").Value = -4379
$ws.Range("H113").Value = 10022.714
$ws.Range("I113").Value = 12057.091
$ws.Range("K113").Value = 12057.091
$ws.Range("M113").Value = -8803.091
$ws.Range("H129").Value = 911.2222
$ws.Range("I129").Value = 1021
$ws.Range("J129").Value = 897.5
$ws.Range("K129").Value = 3063
$ws.Range("L129").Value = 2692.5
$ws.Range("M129").Value = 1937
$ws.Range("N129").Value = -12692.5
$ws.Range("H137").Value = 1560.2222
$ws.Range("I137").Value = 895.3570999999999
$ws.Range("J137").Value = 1983.3182
$ws.Range("K137").Value = 2686.0713
$ws.Range("L137").Value = 5949.9546
$ws.Range("M137").Value = -136.0712999999996
$ws.Range("N137").Value = -11049.9546
$ws.Range("H138").Value = 1809.2142
$ws.Range("J138").Value = 2500
$ws.Range("L138").Value = 7500
$ws.Range("N138").Value = -17780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2831.6765
$ws.Range("I32").Value = 2020.3928
$ws.Range("K32").Value = 2020.3928
$ws.Range("M32").Value = -1733.3928
$ws.Range("H45").Value = 3602683
$ws.Range("I45").Value = 11252547
$ws.Range("J45").Value = 2747
$ws.Range("K45").Value = 11252547
$ws.Range("L45").Value = 2747
$ws.Range("M45").Value = -11252170
$ws.Range("N45").Value = -3501
$ws.Range("H61").Value = 2541.2856
$ws.Range("I61").Value = 1644.1
$ws.Range("K61").Value = 1644.1
$ws.Range("M61").Value = -1432.1
$ws.Range("H74").Value = 1034.0667
$ws.Range("I74").Value = 491.65216
$ws.Range("K74").Value = 491.65216
$ws.Range("M74").Value = 382.34784
$ws.Range("H77").Value = 1034.0667
$ws.Range("I77").Value = 491.65216
$ws.Range("K77").Value = 2458.2608
$ws.Range("M77").Value = 1909.7392
$ws.Range("H122").Value = 1837
$ws.Range("I122").Value = 1012
$ws.Range("J122").Value = 2249.5
$ws.Range("K122").Value = 3036
$ws.Range("L122").Value = 6748.5
$ws.Range("M122").Value = -586
$ws.Range("N122").Value = -11648.5
$ws.Range("H132").Value = 3046.5715
$ws.Range("J132").Value = 4999
$ws.Range("L132").Value = 14997
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 2541.2856
$ws.Range("I136").Value = 1644.1
$ws.Range("K136").Value = 4932.299999999999
$ws.Range("M136").Value = -2382.299999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 183840.45
$ws.Range("I86").Value = 2193.889
$ws.Range("K86").Value = 2193.889
$ws.Range("M86").Value = -1070.889
$ws.Range("H89").Value = 183840.45
$ws.Range("I89").Value = 2193.889
$ws.Range("K89").Value = 10969.445
$ws.Range("M89").Value = -5353.445
$ws.Range("H134").Value = 7146.95
$ws.Range("I134").Value = 7761.1177
$ws.Range("K134").Value = 23283.3531
$ws.Range("M134").Value = -20748.3531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 190
$ws.Range("I7").Value = 190
$ws.Range("K7").Value = 190
$ws.Range("M7").Value = -77
$ws.Range("H31").Value = 1343.1034
$ws.Range("I31").Value = 857.3333
$ws.Range("J31").Value = 1561.7
$ws.Range("K31").Value = 857.3333
$ws.Range("L31").Value = 1561.7
$ws.Range("M31").Value = -562.3333
$ws.Range("N31").Value = -2151.7
$ws.Range("H34").Value = 1343.1034
$ws.Range("I34").Value = 857.3333
$ws.Range("J34").Value = 1561.7
$ws.Range("K34").Value = 857.3333
$ws.Range("L34").Value = 1561.7
$ws.Range("M34").Value = -655.3333
$ws.Range("N34").Value = -1965.7
$ws.Range("H58").Value = 1740318.1
$ws.Range("I58").Value = 3345624.8
$ws.Range("J58").Value = 1236
$ws.Range("K58").Value = 3345624.8
$ws.Range("L58").Value = 1236
$ws.Range("M58").Value = -3345421.8
$ws.Range("N58").Value = -1642
$ws.Range("H99").Value = 3277.6667
$ws.Range("I99").Value = 3277.6667
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3277.6667
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1779.6667
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 3277.6667
$ws.Range("I126").Value = 3277.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9833.000100000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7363.000100000001
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 3648.6155
$ws.Range("I132").Value = 2740
$ws.Range("K132").Value = 8220
$ws.Range("M132").Value = -5690
$ws.Range("H134").Value = 1817.9375
$ws.Range("J134").Value = 2671.8
$ws.Range("L134").Value = 8015.400000000001
$ws.Range("N134").Value = -13085.4
$ws.Range("H136").Value = 1740318.1
$ws.Range("I136").Value = 3345624.8
$ws.Range("J136").Value = 1236
$ws.Range("K136").Value = 10036874.4
$ws.Range("L136").Value = 3708
$ws.Range("M136").Value = -10034324.4
$ws.Range("N136").Value = -8808

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1635.0465
$ws.Range("I68").Value = 755.6667
$ws.Range("K68").Value = 2267.0001
$ws.Range("M68").Value = -1456.0001
$ws.Range("H71").Value = 1635.0465
$ws.Range("I71").Value = 755.6667
$ws.Range("K71").Value = 6801.0003
$ws.Range("M71").Value = -2745.0003
$ws.Range("H131").Value = 14727248
$ws.Range("J131").Value = 24994.104
$ws.Range("L131").Value = 74982.31200000001
$ws.Range("N131").Value = -85062.31200000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7681.684
$ws.Range("J80").Value = 8964.714
$ws.Range("L80").Value = 8964.714
$ws.Range("N80").Value = -10960.714
$ws.Range("H83").Value = 7681.684
$ws.Range("J83").Value = 8964.714
$ws.Range("L83").Value = 44823.57
$ws.Range("N83").Value = -54807.57
$ws.Range("H122").Value = 2412.077
$ws.Range("I122").Value = 2066
$ws.Range("K122").Value = 6198
$ws.Range("M122").Value = -3748
$ws.Range("H132").Value = 4811207
$ws.Range("I132").Value = 7695752
$ws.Range("J132").Value = 3633
$ws.Range("K132").Value = 23087256
$ws.Range("L132").Value = 10899
$ws.Range("M132").Value = -23084726
$ws.Range("N132").Value = -15959

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4004.0588
$ws.Range("J7").Value = 4027.9167
$ws.Range("L7").Value = 4027.9167
$ws.Range("N7").Value = -4251.9167
$ws.Range("H22").Value = 1146.2858
$ws.Range("I22").Value = 894.3333
$ws.Range("J22").Value = 1599.8
$ws.Range("K22").Value = 894.3333
$ws.Range("L22").Value = 1599.8
$ws.Range("M22").Value = -599.3333
$ws.Range("N22").Value = -2189.8
$ws.Range("H27").Value = 1146.2858
$ws.Range("I27").Value = 894.3333
$ws.Range("J27").Value = 1599.8
$ws.Range("K27").Value = 894.3333
$ws.Range("L27").Value = 1599.8
$ws.Range("M27").Value = -787.3333
$ws.Range("N27").Value = -1813.8
$ws.Range("H126").Value = 4004.0588
$ws.Range("J126").Value = 4027.9167
$ws.Range("L126").Value = 12083.7501
$ws.Range("N126").Value = -17023.7501
$ws.Range("H132").Value = 4107.5835
$ws.Range("I132").Value = 1183
$ws.Range("J132").Value = 5082.4443
$ws.Range("K132").Value = 3549
$ws.Range("L132").Value = 15247.3329
$ws.Range("M132").Value = -1019
$ws.Range("N132").Value = -20307.3329
$ws.Range("H136").Value = 2888.4722
$ws.Range("I136").Value = 1792.5927
$ws.Range("K136").Value = 5377.7781
$ws.Range("M136").Value = -2827.7781

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1204.5
$ws.Range("I107").Value = 948.4
$ws.Range("J107").Value = 1844.75
$ws.Range("K107").Value = 2845.2
$ws.Range("L107").Value = 5534.25
$ws.Range("M107").Value = -925.1999999999998
$ws.Range("N107").Value = -9374.25
$ws.Range("H126").Value = 8508.857
$ws.Range("I126").Value = 10260.154
$ws.Range("J126").Value = 5663
$ws.Range("K126").Value = 30780.462
$ws.Range("L126").Value = 16989
$ws.Range("M126").Value = -28310.462
$ws.Range("N126").Value = -21929
$ws.Range("H132").Value = 2152.68
$ws.Range("I132").Value = 1508.75
$ws.Range("J132").Value = 3297.4443
$ws.Range("K132").Value = 4526.25
$ws.Range("L132").Value = 9892.332900000001
$ws.Range("M132").Value = -1996.25
$ws.Range("N132").Value = -14952.3329
$ws.Range("H136").Value = 10103646
$ws.Range("I136").Value = 15434439
$ws.Range("J136").Value = 3194.158
$ws.Range("K136").Value = 46303317
$ws.Range("L136").Value = 9582.474
$ws.Range("M136").Value = -46300767
$ws.Range("N136").Value = -14682.474
